# Updated symbol list (crypto prices / 1h volume change) per upstream diff.
# Each D/E cell holds its value as literal text (matches source inlineStr
# cells), so we force Text number format before assigning to stop Excel's
# COM layer from auto-coercing numeric-looking strings ("312.34", "-0.29%")
# into numbers/percentages, then clear the temporary format so the cell's
# style index is left untouched (matches original unstyled cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "312.34" },
    @{ Cell = "E2"; Value = "-0.29%" },
    @{ Cell = "D3"; Value = "36.88" },
    @{ Cell = "E3"; Value = "-2.14%" },
    @{ Cell = "D4"; Value = "5.118" },
    @{ Cell = "E4"; Value = "-0.28%" },
    @{ Cell = "D5"; Value = "0.07881" },
    @{ Cell = "E5"; Value = "-0.50%" },
    @{ Cell = "D6"; Value = "8.373" },
    @{ Cell = "E6"; Value = "1.07%" },
    @{ Cell = "D7"; Value = "1.870" },
    @{ Cell = "E7"; Value = "-3.26%" },
    @{ Cell = "E8"; Value = "0.98%" },
    @{ Cell = "D9"; Value = "0.9301" },
    @{ Cell = "E9"; Value = "0.88%" },
    @{ Cell = "D10"; Value = "0.1170" },
    @{ Cell = "E10"; Value = "-4.90%" },
    @{ Cell = "D11"; Value = "0.1907" },
    @{ Cell = "E11"; Value = "-1.23%" },
    @{ Cell = "D12"; Value = "0.08881" },
    @{ Cell = "E12"; Value = "-2.87%" },
    @{ Cell = "D13"; Value = "0.03262" },
    @{ Cell = "E13"; Value = "-1.47%" },
    @{ Cell = "D14"; Value = "0.09527" },
    @{ Cell = "E14"; Value = "-1.09%" },
    @{ Cell = "D15"; Value = "0.001380" },
    @{ Cell = "E15"; Value = "0.15%" },
    @{ Cell = "D16"; Value = "0.005814" },
    @{ Cell = "E16"; Value = "1.38%" },
    @{ Cell = "E17"; Value = "-4.67%" },
    @{ Cell = "D18"; Value = "4.388" },
    @{ Cell = "E18"; Value = "-0.53%" },
    @{ Cell = "D19"; Value = "0.3445" },
    @{ Cell = "E19"; Value = "-0.01%" },
    @{ Cell = "D20"; Value = "6.296" },
    @{ Cell = "E20"; Value = "19.89%" },
    @{ Cell = "D21"; Value = "0.1291" },
    @{ Cell = "E21"; Value = "1.48%" },
    @{ Cell = "D22"; Value = "0.2303" },
    @{ Cell = "E22"; Value = "-11.07%" },
    @{ Cell = "D23"; Value = "0.04314" },
    @{ Cell = "E23"; Value = "-1.27%" },
    @{ Cell = "D24"; Value = "0.001193" },
    @{ Cell = "E24"; Value = "-4.49%" },
    @{ Cell = "D25"; Value = "0.004339" },
    @{ Cell = "E25"; Value = "0.64%" },
    @{ Cell = "D26"; Value = "0.0001320" },
    @{ Cell = "E26"; Value = "8.22%" },
    @{ Cell = "D27"; Value = "0.0003950" },
    @{ Cell = "E27"; Value = "-98.12%" },
    @{ Cell = "D39"; Value = "0.02208" },
    @{ Cell = "E39"; Value = "-0.55%" },
    @{ Cell = "D40"; Value = "0.05100" },
    @{ Cell = "E40"; Value = "-0.43%" },
    @{ Cell = "D41"; Value = "0.007561" },
    @{ Cell = "E41"; Value = "1.30%" },
    @{ Cell = "D42"; Value = "0.1369" },
    @{ Cell = "E42"; Value = "0.31%" },
    @{ Cell = "D43"; Value = "0.008311" },
    @{ Cell = "E43"; Value = "-5.32%" },
    @{ Cell = "D44"; Value = "0.001982" },
    @{ Cell = "E44"; Value = "1.21%" },
    @{ Cell = "D45"; Value = "0.007773" },
    @{ Cell = "E45"; Value = "-9.75%" },
    @{ Cell = "D46"; Value = "0.00006304" },
    @{ Cell = "E46"; Value = "-6.28%" },
    @{ Cell = "D47"; Value = "0.00000000744" },
    @{ Cell = "E47"; Value = "-0.72%" },
    @{ Cell = "D48"; Value = "0.002845" },
    @{ Cell = "E48"; Value = "-14.92%" },
    @{ Cell = "D49"; Value = "0.001677" },
    @{ Cell = "E49"; Value = "39.77%" },
    @{ Cell = "D50"; Value = "0.00002084" },
    @{ Cell = "E50"; Value = "-0.72%" },
    @{ Cell = "D51"; Value = "0.0001985" },
    @{ Cell = "E51"; Value = "-0.72%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
